# Apply "Added backup and restore configuration" packet-counter refresh.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "R1" (interface counters)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Columns.Item(1).ColumnWidth = 29.142857142857142

$ws1.Cells.Item(3, 6).Value  = 1954607   # F3 rx_octets
$ws1.Cells.Item(3, 7).Value  = 19690     # G3 rx_unicast
$ws1.Cells.Item(3, 12).Value = 2129092   # L3 tx_octets
$ws1.Cells.Item(3, 13).Value = 17990     # M3 tx_unicast

$ws1.Cells.Item(4, 12).Value = 108914    # L4 tx_octets
$ws1.Cells.Item(4, 13).Value = 964       # M4 tx_unicast

$ws1.Cells.Item(5, 12).Value = 51818     # L5 tx_octets
$ws1.Cells.Item(5, 13).Value = 457       # M5 tx_unicast

$ws1.Cells.Item(6, 12).Value = 48128     # L6 tx_octets
$ws1.Cells.Item(6, 13).Value = 423       # M6 tx_unicast

# ---------------------------------------------------------------------
# Sheet "R3" (interface counters)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("R3")
$ws2.Columns.Item(1).ColumnWidth = 29.142857142857142

$ws2.Cells.Item(3, 6).Value  = 773969    # F3 rx_octets (em0)
$ws2.Cells.Item(3, 12).Value = 1304341   # L3 tx_octets (em0)

$ws2.Cells.Item(4, 12).Value = 2500634   # L4 tx_octets (em1)

$ws2.Cells.Item(5, 6).Value  = 94763     # F5 rx_octets (em2)
$ws2.Cells.Item(5, 12).Value = 640       # L5 tx_octets (em2)

# ---------------------------------------------------------------------
# Sheet "SW1" (interface counters)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SW1")
$ws3.Columns.Item(1).ColumnWidth = 29.142857142857142

$ws3.Cells.Item(3, 2).Value  = 23491     # B3  (GigabitEthernet0/0)
$ws3.Cells.Item(3, 6).Value  = 3205818   # F3
$ws3.Cells.Item(3, 7).Value  = 37281     # G3
$ws3.Cells.Item(3, 12).Value = 2740423   # L3
$ws3.Cells.Item(3, 13).Value = 22444     # M3

$ws3.Cells.Item(4, 12).Value = 2590987   # L4 (GigabitEthernet0/1)
$ws3.Cells.Item(4, 13).Value = 33524     # M4

$ws3.Cells.Item(5, 2).Value  = 32076     # B5 (GigabitEthernet0/2)
$ws3.Cells.Item(5, 5).Value  = 6         # E5
$ws3.Cells.Item(5, 6).Value  = 1882415   # F5
$ws3.Cells.Item(5, 7).Value  = 27689     # G5
$ws3.Cells.Item(5, 12).Value = 566288    # L5
$ws3.Cells.Item(5, 13).Value = 3748      # M5

$ws3.Cells.Item(6, 2).Value  = 14807     # B6 (GigabitEthernet0/3)
$ws3.Cells.Item(6, 5).Value  = 7         # E6
$ws3.Cells.Item(6, 6).Value  = 1007263   # F6
$ws3.Cells.Item(6, 7).Value  = 12525     # G6
$ws3.Cells.Item(6, 12).Value = 1721680   # L6
$ws3.Cells.Item(6, 13).Value = 20985     # M6

$ws3.Cells.Item(7, 12).Value = 300       # L7 (GigabitEthernet1/0)
$ws3.Cells.Item(7, 13).Value = 5         # M7

$ws3.Cells.Item(8, 12).Value = 192       # L8 (GigabitEthernet1/1)
$ws3.Cells.Item(8, 13).Value = 3         # M8

$ws3.Cells.Item(9, 12).Value = 0         # L9 (GigabitEthernet1/2)
$ws3.Cells.Item(9, 13).Value = 0         # M9

$ws3.Cells.Item(10, 12).Value = 0        # L10 (GigabitEthernet1/3)
$ws3.Cells.Item(10, 13).Value = 0        # M10

# New interfaces appended at the bottom of the SW1 sheet.
$newInterfaces = @(
    "GigabitEthernet2/0",
    "GigabitEthernet2/1",
    "GigabitEthernet2/2",
    "GigabitEthernet2/3",
    "GigabitEthernet3/0",
    "GigabitEthernet3/1",
    "GigabitEthernet3/2"
)

$row = 11
foreach ($name in $newInterfaces) {
    $rowRange = $ws3.Range("A" + $row + ":M" + $row)
    $rowRange.HorizontalAlignment = -4108
    $rowRange.VerticalAlignment = -4108

    $ws3.Cells.Item($row, 1).Value  = $name
    $ws3.Cells.Item($row, 2).Value  = 0
    $ws3.Cells.Item($row, 3).Value  = 0
    $ws3.Cells.Item($row, 4).Value  = 0
    $ws3.Cells.Item($row, 5).Value  = 0
    $ws3.Cells.Item($row, 6).Value  = 0
    $ws3.Cells.Item($row, 7).Value  = 0
    $ws3.Cells.Item($row, 8).Value  = -1
    $ws3.Cells.Item($row, 9).Value  = 0
    $ws3.Cells.Item($row, 10).Value = 0
    $ws3.Cells.Item($row, 11).Value = -1
    $ws3.Cells.Item($row, 12).Value = 0
    $ws3.Cells.Item($row, 13).Value = 0

    $row = $row + 1
}
